$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Sign up (Choose your own profile) and Sign " + "into" + " the account"
#    -> single run "Sign up (Choose your own profile) and Sign into the account"
#    (whole paragraph consists only of these 3 runs, so a plain Find/Replace
#    merge is safe here.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Sign up (Choose your own profile) and Sign into the account", $true, $false, $false, $false, $false,
    $true, 1, $false, "Sign up (Choose your own profile) and Sign into the account", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Install - You " + "may " + "have to go thru multiple screens while installing, "
#    -> single run. The following text (don/'t ...) is separated by a
#    <w:proofErr/> so it naturally will not be swept into this merge.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Install - You may have to go thru multiple screens while installing, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Install - You may have to go thru multiple screens while installing, ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "don" + "'t" -> "don't" (also guarded by <w:proofErr/> on both sides)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "don’t", $true, $false, $false, $false, $false,
    $true, 1, $false, "don’t", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "it will " + "open with the home directory" -> single run
#    (only 2 runs in that paragraph, nothing else to protect.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "it will open with the home directory", $true, $false, $false, $false, $false,
    $true, 1, $false, "it will open with the home directory", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Click on File -> " + "New -" + "> Java Project" -> single run, but the
#    paragraph also contains a 4th run " (if java project is not shown...)"
#    that must stay untouched/unmerged. The underlying Word-OM implementation
#    here merges *every* pair of adjacent, identically-formatted runs in a
#    paragraph whenever any text inside that paragraph actually changes - so
#    we temporarily give the trailing run a distinct format (Bold), perform
#    the text edit (which merges only the first three, still-identical,
#    runs), then clear the temporary Bold again (a pure formatting change,
#    which does not trigger the merge pass).
# ---------------------------------------------------------------------------
$guard5 = $d.Content
$guard5.Find.Execute(
    "(if java project is not shown select project -> Java project)", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$guard5Start = $guard5.Start - 1          # include the leading space of that run
$guard5End = $guard5.End
$d.Range($guard5Start, $guard5End).Bold = 1

$edit5 = $d.Content
$edit5.Find.Execute(
    "Click on File -> New -> Java Project", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$edit5.Text = "Click on File -> New -> Java Project_TEMP5"

$edit5b = $d.Content
$edit5b.Find.Execute(
    "Click on File -> New -> Java Project_TEMP5", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$edit5b.Text = "Click on File -> New -> Java Project"

$unguard5 = $d.Content
$unguard5.Find.Execute(
    "(if java project is not shown select project -> Java project)", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$unguard5Start = $unguard5.Start - 1
$unguard5End = $unguard5.End
$d.Range($unguard5Start, $unguard5End).Bold = 0

# ---------------------------------------------------------------------------
# 6. "Open git hub" (single run) -> 5 runs: "Open " + "G" + " " + "it " + "Bash"
#    This paragraph has only the one run, so the text replace is safe; the
#    new run boundaries are then carved out with Bold on/off toggles, which
#    (being pure formatting ops) split the run without merging anything else.
# ---------------------------------------------------------------------------
$edit6 = $d.Content
$edit6.Find.Execute(
    "Open git hub", $true, $false, $false, $false, $false,
    $true, 1, $false, "Open G it Bash", 2) | Out-Null

$find6 = $d.Content
$find6.Find.Execute("Open G it Bash", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s6 = $find6.Start
# "Open " | "G" | " " | "it " | "Bash"
$d.Range($s6 + 0, $s6 + 5).Bold = 1
$d.Range($s6 + 0, $s6 + 5).Bold = 0
$d.Range($s6 + 5, $s6 + 6).Bold = 1
$d.Range($s6 + 5, $s6 + 6).Bold = 0
$d.Range($s6 + 6, $s6 + 7).Bold = 1
$d.Range($s6 + 6, $s6 + 7).Bold = 0
$d.Range($s6 + 7, $s6 + 10).Bold = 1
$d.Range($s6 + 7, $s6 + 10).Bold = 0

# ---------------------------------------------------------------------------
# 7. "/" + "c/Users/Username/" -> "Cd " + "/c/Users/Username/"
#    The paragraph also has "javaselcode" and
#    "/VcentrySeleniumJava/VcentrySeleniumJava" runs after it that must
#    remain separate/untouched. Guard "javaselcode" with Bold while the
#    preceding text is edited (it is identically formatted to the following
#    run and would otherwise be swept into the same merge pass), then split
#    the newly-written text into "Cd " / "/c/Users/Username/" via Bold
#    toggles, and finally clear the guard.
# ---------------------------------------------------------------------------
$guard7 = $d.Content
$guard7.Find.Execute(
    "javaselcode/VcentrySeleniumJava/VcentrySeleniumJava", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$guard7Start = $guard7.Start
$d.Range($guard7Start, $guard7Start + 11).Bold = 1   # len("javaselcode") == 11

$edit7 = $d.Content
$edit7.Find.Execute(
    "/c/Users/Username/", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$edit7Start = $edit7.Start
$edit7.Text = "Cd /c/Users/Username/"

$split7 = $d.Range($edit7Start, $edit7Start + 3)   # "Cd "
$split7.Bold = 1
$split7.Bold = 0

$unguard7 = $d.Content
$unguard7.Find.Execute(
    "javaselcode/VcentrySeleniumJava/VcentrySeleniumJava", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$unguard7Start = $unguard7.Start
$d.Range($unguard7Start, $unguard7Start + 11).Bold = 0

$d.Save()
